$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card9")

# Fix header N1 text (drop trailing space)
$ws.Cells.Item(1, 14).Value = "Correction"

# Copy N1's formatting to new O1 header cell, then set its text
$ws.Cells.Item(1, 14).Copy($ws.Cells.Item(1, 15))
$ws.Cells.Item(1, 15).Value = "Serviced by "

for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
    $ws.Cells.Item($r, 15).Value = "'"
    $ws.Cells.Item($r, 15).ClearFormats()
}
